$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H9").Value = 879
$wsALC.Range("I9").Value = 1325
$wsALC.Range("J9").Value = 433
$wsALC.Range("K9").Value = 1325
$wsALC.Range("L9").Value = 433
$wsALC.Range("M9").Value = -1156
$wsALC.Range("N9").Value = -771
$wsALC.Range("H62").Value = 1427.6666
$wsALC.Range("I62").Value = 1427.6666
$wsALC.Range("K62").Value = 1427.6666
$wsALC.Range("M62").Value = -803.6666
$wsALC.Range("H65").Value = 1427.6666
$wsALC.Range("I65").Value = 1427.6666
$wsALC.Range("K65").Value = 7425
$wsALC.Range("M65").Value = -4018.333000000001
$wsALC.Range("H69").Value = 15750
$wsALC.Range("J69").Value = 15000
$wsALC.Range("L69").Value = 45000
$wsALC.Range("N69").Value = -46748
$wsALC.Range("H72").Value = 15750
$wsALC.Range("J72").Value = 15000
$wsALC.Range("L72").Value = 135000
$wsALC.Range("N72").Value = -143736
$wsALC.Range("H125").Value = 2058
$wsALC.Range("I125").Value = 1516
$wsALC.Range("J125").Value = 2600
$wsALC.Range("K125").Value = 13644
$wsALC.Range("L125").Value = 23400
$wsALC.Range("M125").Value = -11184
$wsALC.Range("N125").Value = -28320
$wsALC.Range("H129").Value = 2336.2
$wsALC.Range("I129").Value = 1268.7142
$wsALC.Range("K129").Value = 3806.1426
$wsALC.Range("M129").Value = 1193.8574
$wsALC.Range("H137").Value = 784.8
$wsALC.Range("I137").Value = 784.8
$wsALC.Range("K137").Value = 2354.4
$wsALC.Range("M137").Value = 195.6000000000004
$wsALC.Range("H138").Value = 1698.375
$wsALC.Range("J138").Value = 4400
$wsALC.Range("L138").Value = 13200
$wsALC.Range("N138").Value = -23480

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H5").Value = 340.16666
$wsARM.Range("J5").Value = 262.5
$wsARM.Range("L5").Value = 262.5
$wsARM.Range("N5").Value = -486.5
$wsARM.Range("H30").Value = 13603
$wsARM.Range("I30").Value = 10404.5
$wsARM.Range("J30").Value = 20000
$wsARM.Range("K30").Value = 10404.5
$wsARM.Range("L30").Value = 20000
$wsARM.Range("M30").Value = -10254.5
$wsARM.Range("N30").Value = -20300
$wsARM.Range("H32").Value = 5162.077
$wsARM.Range("I32").Value = 5162.077
$wsARM.Range("K32").Value = 5162.077
$wsARM.Range("M32").Value = -4875.077
$wsARM.Range("H61").Value = 5659.091
$wsARM.Range("I61").Value = 6064.25
$wsARM.Range("K61").Value = 6064.25
$wsARM.Range("M61").Value = -5852.25
$wsARM.Range("H105").Value = 12345
$wsARM.Range("J105").Value = 12345
$wsARM.Range("L105").Value = 12345
$wsARM.Range("N105").Value = -19333
$wsARM.Range("H110").Value = 5048.269
$wsARM.Range("I110").Value = 5496.933
$wsARM.Range("J110").Value = 4436.4546
$wsARM.Range("K110").Value = 5496.933
$wsARM.Range("L110").Value = 4436.4546
$wsARM.Range("M110").Value = -3451.933
$wsARM.Range("N110").Value = -8526.454600000001
$wsARM.Range("H136").Value = 5659.091
$wsARM.Range("I136").Value = 6064.25
$wsARM.Range("K136").Value = 18192.75
$wsARM.Range("M136").Value = -15642.75
$wsARM.Range("H139").Value = 99999.5
$wsARM.Range("J139").Value = 99999.5
$wsARM.Range("L139").Value = 99999.5
$wsARM.Range("N139").Value = -110279.5

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H4").Value = 340.16666
$wsBSM.Range("J4").Value = 262.5
$wsBSM.Range("L4").Value = 262.5
$wsBSM.Range("N4").Value = -492.5
$wsBSM.Range("H20").Value = 3499.2222
$wsBSM.Range("I20").Value = 3499.2222
$wsBSM.Range("K20").Value = 3499.2222
$wsBSM.Range("M20").Value = -3252.2222

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H22").Value = 2858419.8
$wsCRP.Range("I22").Value = 1317.091
$wsCRP.Range("K22").Value = 1317.091
$wsCRP.Range("M22").Value = -967.0909999999999

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H38").Value = 51.857143
$wsCUL.Range("J38").Value = 17
$wsCUL.Range("L38").Value = 51
$wsCUL.Range("N38").Value = -745
$wsCUL.Range("H107").Value = 144.9
$wsCUL.Range("I107").Value = 136
$wsCUL.Range("J107").Value = 147.125
$wsCUL.Range("K107").Value = 408
$wsCUL.Range("L107").Value = 441.375
$wsCUL.Range("M107").Value = 1512
$wsCUL.Range("N107").Value = -4281.375
$wsCUL.Range("H111").Value = 6600
$wsCUL.Range("I111").Value = 6600
$wsCUL.Range("K111").Value = 19800
$wsCUL.Range("M111").Value = -16733
$wsCUL.Range("H121").Value = 773209.4
$wsCUL.Range("I121").Value = 3247
$wsCUL.Range("J121").Value = 1671498.9
$wsCUL.Range("K121").Value = 9741
$wsCUL.Range("L121").Value = 5014496.699999999
$wsCUL.Range("M121").Value = -8431
$wsCUL.Range("N121").Value = -5017116.699999999
$wsCUL.Range("H131").Value = 1311
$wsCUL.Range("I131").Value = 950
$wsCUL.Range("K131").Value = 2850
$wsCUL.Range("M131").Value = 2190
$wsCUL.Range("H137").Value = 400
$wsCUL.Range("I137").Value = 400
$wsCUL.Range("K137").Value = 1200
$wsCUL.Range("M137").Value = 3900

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H122").Value = 5998
$wsGSM.Range("I122").Value = 5996.5
$wsGSM.Range("J122").Value = 5999.5
$wsGSM.Range("K122").Value = 17989.5
$wsGSM.Range("L122").Value = 17998.5
$wsGSM.Range("M122").Value = -15539.5
$wsGSM.Range("N122").Value = -22898.5

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H22").Value = 7633
$wsLTW.Range("I22").Value = 9999
$wsLTW.Range("J22").Value = 6450
$wsLTW.Range("K22").Value = 9999
$wsLTW.Range("L22").Value = 6450
$wsLTW.Range("M22").Value = -9704
$wsLTW.Range("N22").Value = -7040
$wsLTW.Range("H27").Value = 7633
$wsLTW.Range("I27").Value = 9999
$wsLTW.Range("J27").Value = 6450
$wsLTW.Range("K27").Value = 9999
$wsLTW.Range("L27").Value = 6450
$wsLTW.Range("M27").Value = -9892
$wsLTW.Range("N27").Value = -6664
$wsLTW.Range("H46").Value = 2278
$wsLTW.Range("I46").Value = 1757.8334
$wsLTW.Range("K46").Value = 1757.8334
$wsLTW.Range("M46").Value = -1569.8334
$wsLTW.Range("H105").Value = 0
$wsLTW.Range("J105").Value = 0
$wsLTW.Range("L105").Value = 0
$wsLTW.Range("N105").ClearContents()
$wsLTW.Range("H136").Value = 25999
$wsLTW.Range("I136").Value = 0
$wsLTW.Range("K136").Value = 0
$wsLTW.Range("M136").ClearContents()

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H74").Value = 18710.572
$wsWVR.Range("J74").Value = 17997.25
$wsWVR.Range("L74").Value = 17997.25
$wsWVR.Range("N74").Value = -19869.25
$wsWVR.Range("H77").Value = 18710.572
$wsWVR.Range("J77").Value = 17997.25
$wsWVR.Range("L77").Value = 53991.75
$wsWVR.Range("N77").Value = -63351.75
$wsWVR.Range("H96").Value = 5923.5557
$wsWVR.Range("I96").Value = 5920.3335
$wsWVR.Range("K96").Value = 5920.3335
$wsWVR.Range("M96").Value = -4547.3335
$wsWVR.Range("H107").Value = 450
$wsWVR.Range("I107").Value = 450
$wsWVR.Range("J107").Value = 450
$wsWVR.Range("K107").Value = 1350
$wsWVR.Range("L107").Value = 1350
$wsWVR.Range("M107").Value = 570
$wsWVR.Range("N107").Value = -5190
$wsWVR.Range("H119").Value = 0
$wsWVR.Range("J119").Value = 0
$wsWVR.Range("L119").Value = 0
$wsWVR.Range("N119").ClearContents()
$wsWVR.Range("H122").Value = 2954.4614
$wsWVR.Range("I122").Value = 3078.5557
$wsWVR.Range("J122").Value = 2675.25
$wsWVR.Range("K122").Value = 9235.667099999999
$wsWVR.Range("L122").Value = 8025.75
$wsWVR.Range("M122").Value = -6785.667099999999
$wsWVR.Range("N122").Value = -12925.75
$wsWVR.Range("H126").Value = 2603.45
$wsWVR.Range("I126").Value = 2879.3125
$wsWVR.Range("J126").Value = 1500
$wsWVR.Range("K126").Value = 8637.9375
$wsWVR.Range("L126").Value = 4500
$wsWVR.Range("M126").Value = -6167.9375
$wsWVR.Range("N126").Value = -9440
